# Auto-generated PowerShell Excel COM-interop script
# Applies the URL (pre-signed S3 link) refresh to cells K3, L3, K4, L4
# on the active worksheet, matching the target OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991205821/1749153331348/OwQjboHWXD-camera-video-segment-driver-1749153333848.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSGLGXSZTT%2F20250606%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250606T180107Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIn%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJIMEYCIQDhbQgqpFIZ8BPmbN9gCfgzgElyLh1aIds9TD4wKJzwiQIhAONQhrv9u5F79yTJlVh0Ly35sITt5S760eTBNPm47bJcKt0DCGIQBBoMNzgxMjA0OTQyMjQ0IgxjrVpIddCVvQieBrcqugMA2plH9nZmHg2FJAkoQXJ1A%2BrxfpCbyho0BZCYq2JGp6AnKHHcO61gGDucZrtvirfhBX%2FYOuBSfn06WK6wq6eNvPKF63X7d1AhfE9tN1Zjv5ZgpLLnqO%2BlvGSOWphV0%2Btj7qlY8HbkPTCIPPWxWKZirdBEeL5qqo0uT0URP6xkvfSeFdkzFyjK%2FA72LXhxGxgvt2bj6tWEp6WIqNlrojacFKObL27Qot3HuGvQ1eO0jRQ5DM2%2Fx2NColljXMkxxlDi7poREJUXfE5jr5MIeMoyzSjgVoRkNlznFS%2F3Q1n2URQm0g2CpWjuBcK%2BAx%2FjAv8wxLEKrvXgHK440%2Fukr6NVpTVN739j2E%2FAWzY%2Fe1RFcquX3Ing%2Fdm%2BbwHhgl5JLgudcPRl7bOCnZVOs%2F16rcw4ZdViu3RmZMrps1ogDgs6M7YZpPofvMQiuwBUN71R1ArR7EB4e2YRAT39NcsbceUM9%2B8PKD7IalPKrSzI2WWTGQt29WNnrT2u9S6RyaqSJGbwM6aKPGonfuPevCqUx3Wq%2FIFfGSxeSpZYDKpuUEevIq4G0YzFEAVwvp3BSeIJAmxBAj7iDooRpb9KMOexjMIGOqQB0ZLRoF52Tqun%2FfGaqRio1W1n7n0EgsvutWItQATrBakQHf%2Bc%2FJQp02RHYJ2fXiJRsOJSI%2B4gwtoOWelZsn%2FJCiDMh7KlcJXw4T7ai8IkmFSaGCXVuYbcVtyhXDnM86Df2FXG%2FgUldwcutAKSQyNMqcokrz4qH9KZzelKO1Zd3S7rPKvvlSUttnx1h2HPV07kggp74e%2Be6a3ncT06SdOIBEZSfLg%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2007%20Jun%202025%2002%3A01%3A07%20GMT&X-Amz-Signature=20de6b3fa4bcb286eb31aaec4e6a8b501805c94fe2207fa57ce4df2d3de1bcc9'
$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205821/1749153331348/oV5nVwJUOa-camera-video-segment-1749153333848.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSGLGXSZTT%2F20250606%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250606T180107Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIn%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJIMEYCIQDhbQgqpFIZ8BPmbN9gCfgzgElyLh1aIds9TD4wKJzwiQIhAONQhrv9u5F79yTJlVh0Ly35sITt5S760eTBNPm47bJcKt0DCGIQBBoMNzgxMjA0OTQyMjQ0IgxjrVpIddCVvQieBrcqugMA2plH9nZmHg2FJAkoQXJ1A%2BrxfpCbyho0BZCYq2JGp6AnKHHcO61gGDucZrtvirfhBX%2FYOuBSfn06WK6wq6eNvPKF63X7d1AhfE9tN1Zjv5ZgpLLnqO%2BlvGSOWphV0%2Btj7qlY8HbkPTCIPPWxWKZirdBEeL5qqo0uT0URP6xkvfSeFdkzFyjK%2FA72LXhxGxgvt2bj6tWEp6WIqNlrojacFKObL27Qot3HuGvQ1eO0jRQ5DM2%2Fx2NColljXMkxxlDi7poREJUXfE5jr5MIeMoyzSjgVoRkNlznFS%2F3Q1n2URQm0g2CpWjuBcK%2BAx%2FjAv8wxLEKrvXgHK440%2Fukr6NVpTVN739j2E%2FAWzY%2Fe1RFcquX3Ing%2Fdm%2BbwHhgl5JLgudcPRl7bOCnZVOs%2F16rcw4ZdViu3RmZMrps1ogDgs6M7YZpPofvMQiuwBUN71R1ArR7EB4e2YRAT39NcsbceUM9%2B8PKD7IalPKrSzI2WWTGQt29WNnrT2u9S6RyaqSJGbwM6aKPGonfuPevCqUx3Wq%2FIFfGSxeSpZYDKpuUEevIq4G0YzFEAVwvp3BSeIJAmxBAj7iDooRpb9KMOexjMIGOqQB0ZLRoF52Tqun%2FfGaqRio1W1n7n0EgsvutWItQATrBakQHf%2Bc%2FJQp02RHYJ2fXiJRsOJSI%2B4gwtoOWelZsn%2FJCiDMh7KlcJXw4T7ai8IkmFSaGCXVuYbcVtyhXDnM86Df2FXG%2FgUldwcutAKSQyNMqcokrz4qH9KZzelKO1Zd3S7rPKvvlSUttnx1h2HPV07kggp74e%2Be6a3ncT06SdOIBEZSfLg%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2007%20Jun%202025%2002%3A01%3A07%20GMT&X-Amz-Signature=ead1aaf8a1919b3c505e546d2968fde3ae3c9e0b3c9af7b831b6b5f49ab0f649'
$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991109254/1749147387468/GKgOu1ZeR3-camera-video-segment-driver-1749147389968.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSGLGXSZTT%2F20250606%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250606T180107Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIn%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJIMEYCIQDhbQgqpFIZ8BPmbN9gCfgzgElyLh1aIds9TD4wKJzwiQIhAONQhrv9u5F79yTJlVh0Ly35sITt5S760eTBNPm47bJcKt0DCGIQBBoMNzgxMjA0OTQyMjQ0IgxjrVpIddCVvQieBrcqugMA2plH9nZmHg2FJAkoQXJ1A%2BrxfpCbyho0BZCYq2JGp6AnKHHcO61gGDucZrtvirfhBX%2FYOuBSfn06WK6wq6eNvPKF63X7d1AhfE9tN1Zjv5ZgpLLnqO%2BlvGSOWphV0%2Btj7qlY8HbkPTCIPPWxWKZirdBEeL5qqo0uT0URP6xkvfSeFdkzFyjK%2FA72LXhxGxgvt2bj6tWEp6WIqNlrojacFKObL27Qot3HuGvQ1eO0jRQ5DM2%2Fx2NColljXMkxxlDi7poREJUXfE5jr5MIeMoyzSjgVoRkNlznFS%2F3Q1n2URQm0g2CpWjuBcK%2BAx%2FjAv8wxLEKrvXgHK440%2Fukr6NVpTVN739j2E%2FAWzY%2Fe1RFcquX3Ing%2Fdm%2BbwHhgl5JLgudcPRl7bOCnZVOs%2F16rcw4ZdViu3RmZMrps1ogDgs6M7YZpPofvMQiuwBUN71R1ArR7EB4e2YRAT39NcsbceUM9%2B8PKD7IalPKrSzI2WWTGQt29WNnrT2u9S6RyaqSJGbwM6aKPGonfuPevCqUx3Wq%2FIFfGSxeSpZYDKpuUEevIq4G0YzFEAVwvp3BSeIJAmxBAj7iDooRpb9KMOexjMIGOqQB0ZLRoF52Tqun%2FfGaqRio1W1n7n0EgsvutWItQATrBakQHf%2Bc%2FJQp02RHYJ2fXiJRsOJSI%2B4gwtoOWelZsn%2FJCiDMh7KlcJXw4T7ai8IkmFSaGCXVuYbcVtyhXDnM86Df2FXG%2FgUldwcutAKSQyNMqcokrz4qH9KZzelKO1Zd3S7rPKvvlSUttnx1h2HPV07kggp74e%2Be6a3ncT06SdOIBEZSfLg%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2007%20Jun%202025%2002%3A01%3A07%20GMT&X-Amz-Signature=4a77aa81505a30eb3e24d40fc84307fe2f0783b5e767913dd149e960e39e28ed'
$ws.Range("L4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991109254/1749147387468/dQoeD1EUOr-camera-video-segment-1749147389968.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSGLGXSZTT%2F20250606%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250606T180107Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEIn%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJIMEYCIQDhbQgqpFIZ8BPmbN9gCfgzgElyLh1aIds9TD4wKJzwiQIhAONQhrv9u5F79yTJlVh0Ly35sITt5S760eTBNPm47bJcKt0DCGIQBBoMNzgxMjA0OTQyMjQ0IgxjrVpIddCVvQieBrcqugMA2plH9nZmHg2FJAkoQXJ1A%2BrxfpCbyho0BZCYq2JGp6AnKHHcO61gGDucZrtvirfhBX%2FYOuBSfn06WK6wq6eNvPKF63X7d1AhfE9tN1Zjv5ZgpLLnqO%2BlvGSOWphV0%2Btj7qlY8HbkPTCIPPWxWKZirdBEeL5qqo0uT0URP6xkvfSeFdkzFyjK%2FA72LXhxGxgvt2bj6tWEp6WIqNlrojacFKObL27Qot3HuGvQ1eO0jRQ5DM2%2Fx2NColljXMkxxlDi7poREJUXfE5jr5MIeMoyzSjgVoRkNlznFS%2F3Q1n2URQm0g2CpWjuBcK%2BAx%2FjAv8wxLEKrvXgHK440%2Fukr6NVpTVN739j2E%2FAWzY%2Fe1RFcquX3Ing%2Fdm%2BbwHhgl5JLgudcPRl7bOCnZVOs%2F16rcw4ZdViu3RmZMrps1ogDgs6M7YZpPofvMQiuwBUN71R1ArR7EB4e2YRAT39NcsbceUM9%2B8PKD7IalPKrSzI2WWTGQt29WNnrT2u9S6RyaqSJGbwM6aKPGonfuPevCqUx3Wq%2FIFfGSxeSpZYDKpuUEevIq4G0YzFEAVwvp3BSeIJAmxBAj7iDooRpb9KMOexjMIGOqQB0ZLRoF52Tqun%2FfGaqRio1W1n7n0EgsvutWItQATrBakQHf%2Bc%2FJQp02RHYJ2fXiJRsOJSI%2B4gwtoOWelZsn%2FJCiDMh7KlcJXw4T7ai8IkmFSaGCXVuYbcVtyhXDnM86Df2FXG%2FgUldwcutAKSQyNMqcokrz4qH9KZzelKO1Zd3S7rPKvvlSUttnx1h2HPV07kggp74e%2Be6a3ncT06SdOIBEZSfLg%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2007%20Jun%202025%2002%3A01%3A07%20GMT&X-Amz-Signature=c841d84d0834dd3c7109607704994f3ab424806a0e6a451d2a936dfb339a584e'
